$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.461.21'
$ws.Range('E2').Value = '  +4.18%  '
$ws.Range('D3').Value = '3.499.63'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '590.58'
$ws.Range('E5').Value = '  +3.62%  '
$ws.Range('D6').Value = '170.01'
$ws.Range('E6').Value = '  +8.11%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.501.57'
$ws.Range('E8').Value = '  +2.33%  '
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').Value = '0.124'
$ws.Range('E11').Value = '  +4.08%  '
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('D13').Value = '4.103.37'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = '28.03'
$ws.Range('E15').Value = '  +3.18%  '
$ws.Range('D16').Value = '66.463.35'
$ws.Range('E16').Value = '  +4.02%  '
$ws.Range('D17').Value = '0.0000177'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '3.616.26'
$ws.Range('E18').Value = '  +6.00%  '
$ws.Range('D19').Value = '6.28'
$ws.Range('E19').Value = '  +2.64%  '
$ws.Range('D20').Value = '14.03'
$ws.Range('E20').Value = '  +2.80%  '
$ws.Range('D21').Value = '388.04'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').Value = '7.98'
$ws.Range('E22').Value = '  +2.67%  '
$ws.Range('E23').Value = '  +2.46%  '
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '0.528'
$ws.Range('E25').Value = '  +1.72%  '
$ws.Range('E26').Value = '  +6.72%  '
$ws.Range('D27').Value = '10.06'
$ws.Range('E27').Value = '  +3.74%  '
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '6.39'
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('E31').Value = '  +5.79%  '
$ws.Range('E32').Value = '  +3.46%  '
$ws.Range('D33').Value = '23.47'
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('E34').Value = '  +5.68%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').Value = '160.87'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').Value = '0.903'
$ws.Range('E38').Value = '  +8.36%  '
$ws.Range('D39').Value = '1.94'
$ws.Range('E39').Value = '  +5.33%  '
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('D41').Value = '26.50'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '6.70'
$ws.Range('E42').Value = '  +4.61%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '27.15'
$ws.Range('E43').Value = '  +5.26%  '
$ws.Range('D44').Value = '4.60'
$ws.Range('E44').Value = '  +4.11%  '
$ws.Range('D45').Value = '2.809.27'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '43.48'
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('E47').Value = '  +3.01%  '
$ws.Range('D48').Value = '356.43'
$ws.Range('E48').Value = '  +8.60%  '
$ws.Range('D49').Value = '2.53'
$ws.Range('E49').Value = '  +8.55%  '
$ws.Range('E50').Value = '  +6.04%  '
$ws.Range('D51').Value = '32.80'
$ws.Range('E51').Value = '  +9.05%  '
